# Update the "Estado de Cuenta" sheet: swap the two "Periodo Mora" / "Valor Mora"
# rows so that period 2201 (with its corresponding overdue value) now appears
# before period 2112, matching the refreshed source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16 (first worker period/value) and Row 17 (second worker period/value)
# currently hold period "2112"/18170 and "2201"/36341 respectively.
# The data refresh swaps them: row 16 -> "2201"/36341, row 17 -> "2112"/18170.
$ws.Range("E16").Value = "2201"
$ws.Range("F16").Value = 36341

$ws.Range("E17").Value = "2112"
$ws.Range("F17").Value = 18170
